$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2022" column (M) that mirrors the existing "2021" column (L):
# first copy L3:L11's formatting into M3:M11, then overwrite the values in
# the new column with the 2022 figures (same figures as 2021 for every
# metric except the year header itself).
$xlPasteFormats = -4122
$ws.Range("L3:L11").Copy()
$ws.Range("M3:M11").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("M4").Value2 = 2022
$ws.Range("M5").Value2 = 0.86
$ws.Range("M6").Value2 = 1.07
$ws.Range("M7").Value2 = 25.27
$ws.Range("M8").Value2 = 14
$ws.Range("M9").Value2 = 0.12
$ws.Range("M10").Value2 = 21.74
$ws.Range("M11").Value2 = 9.4600000000000009

# Move the active selection, matching the author's recorded cursor position.
$ws.Range("N6").Select() | Out-Null
